$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.807599666666667
$ws.Range("H2").Value = 5.422799
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 40.70766766666667
$ws.Range("N2").Value = 122.123003
$ws.Range("O2").Value = 0.3776398983502007
$ws.Range("P2").Value = 0.3776398983502007
$ws.Range("Q2").Value = 73.5831665050441
$ws.Range("R2").Value = 662.2484985453971
$ws.Range("S2").Value = 0.3776398983502007
$ws.Range("T2").Value = 0.3776398983502007

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.807599666666667
$ws.Range("H3").Value = 5.422799
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 39.715023
$ws.Range("N3").Value = 119.145069
$ws.Range("O3").Value = 0.3684312589831062
$ws.Range("P3").Value = 0.3684312589831062
$ws.Range("Q3").Value = 71.788862336459
$ws.Range("R3").Value = 646.0997610281311
$ws.Range("S3").Value = 0.3684312589831062
$ws.Range("T3").Value = 0.3684312589831062

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.807599666666667
$ws.Range("H4").Value = 5.422799
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 27.37224266666666
$ws.Range("N4").Value = 82.11672799999999
$ws.Range("O4").Value = 0.253928842666693
$ws.Range("P4").Value = 0.253928842666693
$ws.Range("Q4").Value = 49.47805672018578
$ws.Range("R4").Value = 445.302510481672
$ws.Range("S4").Value = 0.253928842666693
$ws.Range("T4").Value = 0.253928842666693
